# Added ifoCAST full series evaluation: the naive QoQ error "staircase"
# series gains one more diagonal of data. Every existing row's values shift
# one column to the left (B2->dropped, C2->B2, D2->C2, ...), rows that reach
# a new data point get a freshly computed value appended in their last
# column, and rows that have no new value simply shrink by one cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.2000313112857
$ws.Range("C2").Value = -7.713854356263806
$ws.Range("D2").Value = 0.4730551773759682
$ws.Range("E2").Value = 1.112750016219461
$ws.Range("F2").Value = -0.688116008571576
$ws.Range("G2").Value = -1.647050705078689
$ws.Range("H2").Value = 1.1866675198023
$ws.Range("I2").Value = -0.4824371439038641
$ws.Range("J2").Value = 0.5377352938393787
$ws.Range("K2").Value = -1.176054163059933

$ws.Range("B3").Value = -8.028779524684778
$ws.Range("C3").Value = 0.1581300089549967
$ws.Range("D3").Value = 0.7978248477984897
$ws.Range("E3").Value = -1.003041176992548
$ws.Range("F3").Value = -1.96197587349966
$ws.Range("G3").Value = 0.8717423513813284
$ws.Range("H3").Value = -0.7973623123248356
$ws.Range("I3").Value = 0.2228101254184071
$ws.Range("J3").Value = -1.490979331480904
$ws.Range("K3").Value = 1.337376042118068

$ws.Range("B4").Value = -0.898244951227533
$ws.Range("C4").Value = -0.2585501123840399
$ws.Range("D4").Value = -2.059416137175077
$ws.Range("E4").Value = -3.01835083368219
$ws.Range("F4").Value = -0.1846326088012012
$ws.Range("G4").Value = -1.853737272507365
$ws.Range("H4").Value = -0.8335648347641225
$ws.Range("I4").Value = -2.547354291663434
$ws.Range("J4").Value = 0.281001081935538
$ws.Range("K4").Value = -1.954331517930581

$ws.Range("B5").Value = 0.6120962635868906
$ws.Range("C5").Value = -1.188769761204147
$ws.Range("D5").Value = -2.14770445771126
$ws.Range("E5").Value = 0.6860137671697293
$ws.Range("F5").Value = -0.9830908965364347
$ws.Range("G5").Value = 0.03708154120680809
$ws.Range("H5").Value = -1.676707915692503
$ws.Range("I5").Value = 1.151647457906469
$ws.Range("J5").Value = -1.083685141959651
$ws.Range("K5").Value = 0.2676499159580117

$ws.Range("B6").Value = -1.19942097753403
$ws.Range("C6").Value = -2.158355674041143
$ws.Range("D6").Value = 0.6753625508398458
$ws.Range("E6").Value = -0.9937421128663182
$ws.Range("F6").Value = 0.02643032487692459
$ws.Range("G6").Value = -1.687359132022387
$ws.Range("H6").Value = 1.140996241576585
$ws.Range("I6").Value = -1.094336358289534
$ws.Range("J6").Value = 0.2569986996281282
$ws.Range("K6").Value = -0.3644392301887736

$ws.Range("B7").Value = -2.266970206257284
$ws.Range("C7").Value = 0.5667480186237051
$ws.Range("D7").Value = -1.102356645082459
$ws.Range("E7").Value = -0.08218420733921622
$ws.Range("F7").Value = -1.795973664238528
$ws.Range("G7").Value = 1.032381709360444
$ws.Range("H7").Value = -1.202950890505675
$ws.Range("I7").Value = 0.1483841674119874
$ws.Range("J7").Value = -0.4730537624049144
$ws.Range("K7").Value = 0.02750693478591659

$ws.Range("B8").Value = 0.774849739591444
$ws.Range("C8").Value = -0.89425492411472
$ws.Range("D8").Value = 0.1259175136285228
$ws.Range("E8").Value = -1.587871943270789
$ws.Range("F8").Value = 1.240483430328183
$ws.Range("G8").Value = -0.994849169537936
$ws.Range("H8").Value = 0.3564858883797264
$ws.Range("I8").Value = -0.2649520414371754
$ws.Range("J8").Value = 0.2356086557536556
$ws.Range("K8").Value = -0.4322994165924858

$ws.Range("B9").Value = -0.7777567537409195
$ws.Range("C9").Value = 0.2424156840023232
$ws.Range("D9").Value = -1.471373772896988
$ws.Range("E9").Value = 1.356981600701984
$ws.Range("F9").Value = -0.8783509991641355
$ws.Range("G9").Value = 0.4729840587535268
$ws.Range("H9").Value = -0.148453871063375
$ws.Range("I9").Value = 0.352106826127456
$ws.Range("J9").Value = -0.3158012462186854
$ws.Range("K9").Value = 0.1131997290193177

$ws.Range("B10").Value = 0.2879090979994584
$ws.Range("C10").Value = -1.425880358899853
$ws.Range("D10").Value = 1.402475014699119
$ws.Range("E10").Value = -0.8328575851670005
$ws.Range("F10").Value = 0.5184774727506619
$ws.Range("G10").Value = -0.1029604570662399
$ws.Range("H10").Value = 0.3976002401245912
$ws.Range("I10").Value = -0.2703078322215502
$ws.Range("J10").Value = 0.1586931430164528
$ws.Range("K10").Value = 0.2163646915946629

$ws.Range("B11").Value = -1.443434480259818
$ws.Range("C11").Value = 1.384920893339154
$ws.Range("D11").Value = -0.8504117065269649
$ws.Range("E11").Value = 0.5009233513906975
$ws.Range("F11").Value = -0.1205145784262043
$ws.Range("G11").Value = 0.3800461187646267
$ws.Range("H11").Value = -0.2878619535815147
$ws.Range("I11").Value = 0.1411390216564884
$ws.Range("J11").Value = 0.1988105702346985
$ws.Range("K11").Value = 0.322776941072984

$ws.Range("B12").Value = 1.473028212290161
$ws.Range("C12").Value = -0.7623043875759586
$ws.Range("D12").Value = 0.5890306703417038
$ws.Range("E12").Value = -0.0324072594751981
$ws.Range("F12").Value = 0.4681534377156329
$ws.Range("G12").Value = -0.1997546346305085
$ws.Range("H12").Value = 0.2292463406074946
$ws.Range("I12").Value = 0.2869178891857047
$ws.Range("J12").Value = 0.4108842600239903
$ws.Range("K12").Value = -0.4825338632108016

$ws.Range("B13").Value = -0.5837297540881751
$ws.Range("C13").Value = 0.7676053038294873
$ws.Range("D13").Value = 0.1461673740125855
$ws.Range("E13").Value = 0.6467280712034165
$ws.Range("F13").Value = -0.02118000114272489
$ws.Range("G13").Value = 0.4078209740952782
$ws.Range("H13").Value = 0.4654925226734883
$ws.Range("I13").Value = 0.5894588935117738
$ws.Range("J13").Value = -0.303959229723018
$ws.Range("K13").Value = 0.4661714972207444

$ws.Range("B14").Value = 1.679632531582137
$ws.Range("C14").Value = 1.058194601765235
$ws.Range("D14").Value = 1.558755298956066
$ws.Range("E14").Value = 0.8908472266099251
$ws.Range("F14").Value = 1.319848201847928
$ws.Range("G14").Value = 1.377519750426138
$ws.Range("H14").Value = 1.501486121264424
$ws.Range("I14").Value = 0.608067998029632
$ws.Range("J14").Value = 1.378198724973394
$ws.Range("K14").Value = 1.11229800409388

$ws.Range("B15").Value = 0.1142203657994787
$ws.Range("C15").Value = 0.6147810629903097
$ws.Range("D15").Value = -0.0531270093558317
$ws.Range("E15").Value = 0.3758739658821714
$ws.Range("F15").Value = 0.4335455144603815
$ws.Range("G15").Value = 0.557511885298667
$ws.Range("H15").Value = -0.3359062379361248
$ws.Range("I15").Value = 0.4342244890076376
$ws.Range("J15").Value = 0.1683237681281231
$ws.Range("K15").ClearContents()

$ws.Range("B16").Value = 0.6187489605034189
$ws.Range("C16").Value = -0.04915911184272259
$ws.Range("D16").Value = 0.3798418633952805
$ws.Range("E16").Value = 0.4375134119734906
$ws.Range("F16").Value = 0.5614797828117761
$ws.Range("G16").Value = -0.3319383404230157
$ws.Range("H16").Value = 0.4381923865207467
$ws.Range("I16").Value = 0.1722916656412322
$ws.Range("J16").ClearContents()

$ws.Range("B17").Value = 0.08648097832751878
$ws.Range("C17").Value = 0.5154819535655218
$ws.Range("D17").Value = 0.573153502143732
$ws.Range("E17").Value = 0.6971198729820175
$ws.Range("F17").Value = -0.1962982502527744
$ws.Range("G17").Value = 0.5738324766909881
$ws.Range("H17").Value = 0.3079317558114735
$ws.Range("I17").ClearContents()

$ws.Range("B18").Value = 0.2746757717098572
$ws.Range("C18").Value = 0.3323473202880673
$ws.Range("D18").Value = 0.4563136911263528
$ws.Range("E18").Value = -0.4371044321084391
$ws.Range("F18").Value = 0.3330262948353234
$ws.Range("G18").Value = 0.06712557395580883
$ws.Range("H18").ClearContents()

$ws.Range("B19").Value = 0.2870161050359709
$ws.Range("C19").Value = 0.4109824758742565
$ws.Range("D19").Value = -0.4824356473605354
$ws.Range("E19").Value = 0.287695079583227
$ws.Range("F19").Value = 0.02179435870371246
$ws.Range("G19").ClearContents()

$ws.Range("B20").Value = 0.3441210539382026
$ws.Range("C20").Value = -0.5492970692965893
$ws.Range("D20").Value = 0.2208336576471732
$ws.Range("E20").Value = -0.04506706323234141
$ws.Range("F20").ClearContents()

$ws.Range("B21").Value = -0.5788832716533059
$ws.Range("C21").Value = 0.1912474552904566
$ws.Range("D21").Value = -0.07465326558905801
$ws.Range("E21").ClearContents()

$ws.Range("B22").Value = 0.1730967985608157
$ws.Range("C22").Value = -0.0928039223186989
$ws.Range("D22").ClearContents()

$ws.Range("B23").Value = -0.1108357465673982
$ws.Range("C23").ClearContents()

$ws.Range("B24").ClearContents()

